# Update the cryptocurrency price (column D) and volume-change (column E)
# figures on the active sheet to reflect the latest scrape.
#
# Note: several "Price" values (column D) are plain decimal numbers
# (e.g. 401.98, 0.0361, 1.00). If assigned as-is, Excel would
# auto-convert them to numeric values and normalize/clip their textual
# representation (e.g. "1.00" -> 1, "0.0361" -> 3.61E-02), which would
# not match the original text-formatted cell contents. Prefixing the
# value with a leading single quote forces Excel to store it verbatim
# as text while keeping the visible text itself free of the quote.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.507.46'
$ws.Range("E2").Value = '  +9.65%  '
$ws.Range("D3").Value = '3.239.44'
$ws.Range("E3").Value = '  +4.51%  '
$ws.Range("D5").Value = '''401.98'
$ws.Range("E5").Value = '  +4.56%  '
$ws.Range("D6").Value = '''112.05'
$ws.Range("E6").Value = '  +8.81%  '
$ws.Range("D7").Value = '''0.559'
$ws.Range("E7").Value = '  +3.43%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  +7.49%  '
$ws.Range("D10").Value = '''39.75'
$ws.Range("E10").Value = '  +7.67%  '
$ws.Range("D11").Value = '''0.0904'
$ws.Range("E11").Value = '  +5.62%  '
$ws.Range("E12").Value = '  +2.23%  '
$ws.Range("D13").Value = '3.747.03'
$ws.Range("E13").Value = '  +4.55%  '
$ws.Range("D14").Value = '''19.26'
$ws.Range("E14").Value = '  +3.36%  '
$ws.Range("D15").Value = '''8.13'
$ws.Range("E15").Value = '  +3.62%  '
$ws.Range("E16").Value = '  +7.59%  '
$ws.Range("D17").Value = '3.240.43'
$ws.Range("E17").Value = '  +4.79%  '
$ws.Range("D18").Value = '''10.67'
$ws.Range("E18").Value = '  -4.37%  '
$ws.Range("D19").Value = '56.337.02'
$ws.Range("E19").Value = '  +9.32%  '
$ws.Range("D20").Value = '''3.46'
$ws.Range("E20").Value = '  +3.81%  '
$ws.Range("E21").Value = '  +7.85%  '
$ws.Range("D22").Value = '''13.20'
$ws.Range("E22").Value = '  +6.66%  '
$ws.Range("D23").Value = '''310.75'
$ws.Range("E23").Value = '  +16.83%  '
$ws.Range("D24").Value = '''75.06'
$ws.Range("E24").Value = '  +7.30%  '
$ws.Range("D25").Value = '''3.27'
$ws.Range("E25").Value = '  +4.47%  '
$ws.Range("D26").Value = '''8.28'
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("D27").Value = '''28.57'
$ws.Range("E27").Value = '  +5.61%  '
$ws.Range("E28").Value = '  +3.35%  '
$ws.Range("E29").Value = '  +2.88%  '
$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("E31").Value = '  +5.22%  '
$ws.Range("D32").Value = '''11.44'
$ws.Range("E32").Value = '  +10.77%  '
$ws.Range("E33").Value = '  +6.39%  '
$ws.Range("D34").Value = '''36.83'
$ws.Range("E34").Value = '  +4.18%  '
$ws.Range("E35").Value = '  +2.21%  '
$ws.Range("D36").Value = '''51.49'
$ws.Range("E36").Value = '  +2.39%  '
$ws.Range("D37").Value = '''3.57'
$ws.Range("E37").Value = '  +6.04%  '
$ws.Range("D38").Value = '''3.12'
$ws.Range("E38").Value = '  +23.89%  '
$ws.Range("D39").Value = '''0.999'
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '''135.12'
$ws.Range("E40").Value = '  +4.90%  '
$ws.Range("E41").Value = '  +2.93%  '
$ws.Range("D42").Value = '''4.06'
$ws.Range("E42").Value = '  +11.10%  '
$ws.Range("D43").Value = '''17.27'
$ws.Range("E43").Value = '  +4.21%  '
$ws.Range("D44").Value = '''0.120'
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("D45").Value = '''0.284'
$ws.Range("E45").Value = '  -4.74%  '
$ws.Range("D46").Value = '''22.71'
$ws.Range("E46").Value = '  +1.17%  '
$ws.Range("D47").Value = '''2.12'
$ws.Range("E47").Value = '  +40.16%  '
$ws.Range("D48").Value = '2.161.78'
$ws.Range("E48").Value = '  +5.12%  '
$ws.Range("E49").Value = '  +0.82%  '
$ws.Range("D50").Value = '''2.42'
$ws.Range("E50").Value = '  -1.72%  '
$ws.Range("D51").Value = '''0.0361'
$ws.Range("E51").Value = '  +9.72%  '
